# Rewrite the "type" (column B) and "value" (column C) columns of the
# Sports output sheet so that:
#   - column B encodes both the original club/uil prefix and the
#     coed/boys/girls suffix derived from column C, e.g.
#       club-sports + Tennis-Coed  -> sports_club_coed
#       uil-sports  + Tennis-Boys -> sports_uil_boys
#   - column C is collapsed down to just the sport name ("Tennis"),
#     since the boys/girls/coed qualifier now lives in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $typeCell  = $ws.Cells.Item($r, 2)
    $valueCell = $ws.Cells.Item($r, 3)

    $oldType  = [string]$typeCell.Value2
    $oldValue = [string]$valueCell.Value2

    if ([string]::IsNullOrEmpty($oldType) -and [string]::IsNullOrEmpty($oldValue)) {
        continue
    }

    if ($oldType -eq "club-sports") {
        $prefix = "club"
    } elseif ($oldType -eq "uil-sports") {
        $prefix = "uil"
    } else {
        $prefix = $oldType
    }

    if ($oldValue -like "*-Coed") {
        $suffix = "coed"
    } elseif ($oldValue -like "*-Boys") {
        $suffix = "boys"
    } elseif ($oldValue -like "*-Girls") {
        $suffix = "girls"
    } else {
        $suffix = $null
    }

    # Base sport name with the qualifier suffix stripped off.
    $sportName = $oldValue
    foreach ($q in @("-Coed", "-Boys", "-Girls")) {
        if ($sportName -like "*$q") {
            $sportName = $sportName.Substring(0, $sportName.Length - $q.Length)
        }
    }

    if ($suffix) {
        $newType = "sports_{0}_{1}" -f $prefix, $suffix
    } else {
        $newType = $oldType
    }

    $typeCell.Value = $newType
    $valueCell.Value = $sportName
}
